$wb = $excel.ActiveWorkbook

# --- RR_GuiContent: column F ("Deprecated") values change from the
# boolean FALSE to the text "no" for every data row (F2:F48). Setting
# the whole range at once also creates/reuses the shared string "no". ---
$wsGuiContent = $wb.Worksheets.Item("RR_GuiContent")
$wsGuiContent.Range("F2:F48").Value = "no"

# --- RR_Resources was previously the active/selected tab; touch it so
# its prior selection (I9) is preserved, then hand activation over to
# RR_GuiContent, which becomes the new active tab/selected sheet. ---
$wsResources = $wb.Worksheets.Item("RR_Resources")
$wsResources.Activate()
$wsResources.Range("I9").Select()

# --- RR_GuiContent becomes the active sheet, with F2:F48 selected
# (active cell F2). This also updates workbook.xml's activeTab and
# clears tabSelected from RR_Resources' sheetView. ---
$wsGuiContent.Activate()
$wsGuiContent.Range("F2:F48").Select()
